# Auto-generated edit script: update market-price derived values per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 4712.304
$ws.Range("I41").Value = 378.64285
$ws.Range("J41").Value = 11453.556
$ws.Range("K41").Value = 378.64285
$ws.Range("L41").Value = 11453.556
$ws.Range("M41").Value = 61.35714999999999
$ws.Range("N41").Value = -12333.556

$ws.Range("H64").Value = 3828.7666
$ws.Range("I64").Value = 3824.8333
$ws.Range("J64").Value = 3831.389
$ws.Range("K64").Value = 3824.8333
$ws.Range("L64").Value = 3831.389
$ws.Range("M64").Value = -3576.8333
$ws.Range("N64").Value = -4327.389

$ws.Range("H67").Value = 3828.7666
$ws.Range("I67").Value = 3824.8333
$ws.Range("J67").Value = 3831.389
$ws.Range("K67").Value = 3824.8333
$ws.Range("L67").Value = 3831.389
$ws.Range("M67").Value = -2966.8333
$ws.Range("N67").Value = -5547.389

$ws.Range("H76").Value = 8553.375
$ws.Range("I76").Value = 9093.736999999999
$ws.Range("K76").Value = 9093.736999999999
$ws.Range("M76").Value = -8778.736999999999

$ws.Range("H79").Value = 8553.375
$ws.Range("I79").Value = 9093.736999999999
$ws.Range("K79").Value = 9093.736999999999
$ws.Range("M79").Value = -8001.736999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13739.168
$ws.Range("I32").Value = 6506.6987
$ws.Range("J32").Value = 37737.816
$ws.Range("K32").Value = 6506.6987
$ws.Range("L32").Value = 37737.816
$ws.Range("M32").Value = -6219.6987
$ws.Range("N32").Value = -38311.816

$ws.Range("H45").Value = 1524.1111
$ws.Range("I45").Value = 1348.7693
$ws.Range("J45").Value = 1980
$ws.Range("K45").Value = 1348.7693
$ws.Range("L45").Value = 1980
$ws.Range("M45").Value = -971.7692999999999
$ws.Range("N45").Value = -2734

$ws.Range("H63").Value = 7354.4614
$ws.Range("I63").Value = 8867.556
$ws.Range("J63").Value = 3950
$ws.Range("K63").Value = 8867.556
$ws.Range("L63").Value = 3950
$ws.Range("M63").Value = -8181.556
$ws.Range("N63").Value = -5322

$ws.Range("H66").Value = 7354.4614
$ws.Range("I66").Value = 8867.556
$ws.Range("J66").Value = 3950
$ws.Range("K66").Value = 44337.78
$ws.Range("L66").Value = 19750
$ws.Range("M66").Value = -40905.78
$ws.Range("N66").Value = -26614

$ws.Range("H88").Value = 1897.9348
$ws.Range("I88").Value = 1920.9722
$ws.Range("J88").Value = 1815
$ws.Range("K88").Value = 1920.9722
$ws.Range("L88").Value = 1815
$ws.Range("M88").Value = -1514.9722
$ws.Range("N88").Value = -2627

$ws.Range("H91").Value = 1897.9348
$ws.Range("I91").Value = 1920.9722
$ws.Range("J91").Value = 1815
$ws.Range("K91").Value = 1920.9722
$ws.Range("L91").Value = 1815
$ws.Range("M91").Value = -516.9721999999999
$ws.Range("N91").Value = -4623

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2436.3833
$ws.Range("I105").Value = 2171.44
$ws.Range("K105").Value = 2171.44
$ws.Range("M105").Value = -424.4400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 84395.25
$ws.Range("I16").Value = 111815.89
$ws.Range("J16").Value = 2133.3333
$ws.Range("K16").Value = 111815.89
$ws.Range("L16").Value = 2133.3333
$ws.Range("M16").Value = -111528.89
$ws.Range("N16").Value = -2707.3333

$ws.Range("H62").Value = 5061.5557
$ws.Range("J62").Value = 5319.25
$ws.Range("L62").Value = 5319.25
$ws.Range("N62").Value = -6567.25

$ws.Range("H65").Value = 5061.5557
$ws.Range("J65").Value = 5319.25
$ws.Range("L65").Value = 26596.25
$ws.Range("N65").Value = -32836.25

$ws.Range("H113").Value = 84395.25
$ws.Range("I113").Value = 111815.89
$ws.Range("J113").Value = 2133.3333
$ws.Range("K113").Value = 111815.89
$ws.Range("L113").Value = 2133.3333
$ws.Range("M113").Value = -109645.89
$ws.Range("N113").Value = -6473.3333

$ws.Range("H134").Value = 21741336
$ws.Range("I134").Value = 50001496
$ws.Range("K134").Value = 150004488
$ws.Range("M134").Value = -150001953

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 204.6
$ws.Range("I8").Value = 204.6
$ws.Range("K8").Value = 613.8
$ws.Range("M8").Value = -474.8

$ws.Range("H38").Value = 52.36842
$ws.Range("J38").Value = 43.166668
$ws.Range("L38").Value = 129.500004
$ws.Range("N38").Value = -823.500004

$ws.Range("H75").Value = 1743.3334
$ws.Range("J75").Value = 1743.3334
$ws.Range("L75").Value = 5230.0002
$ws.Range("N75").Value = -7226.0002

$ws.Range("H78").Value = 1743.3334
$ws.Range("J78").Value = 1743.3334
$ws.Range("L78").Value = 15690.0006
$ws.Range("N78").Value = -25674.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6529.4893
$ws.Range("I70").Value = 3845.2122
$ws.Range("J70").Value = 12856.714
$ws.Range("K70").Value = 3845.2122
$ws.Range("L70").Value = 12856.714
$ws.Range("M70").Value = -3575.2122
$ws.Range("N70").Value = -13396.714

$ws.Range("H73").Value = 6529.4893
$ws.Range("I73").Value = 3845.2122
$ws.Range("J73").Value = 12856.714
$ws.Range("K73").Value = 3845.2122
$ws.Range("L73").Value = 12856.714
$ws.Range("M73").Value = -2909.2122
$ws.Range("N73").Value = -14728.714

$ws.Range("H80").Value = 4758.75
$ws.Range("I80").Value = 8522.5
$ws.Range("J80").Value = 2876.875
$ws.Range("K80").Value = 8522.5
$ws.Range("L80").Value = 2876.875
$ws.Range("M80").Value = -7524.5
$ws.Range("N80").Value = -4872.875

$ws.Range("H83").Value = 4758.75
$ws.Range("I83").Value = 8522.5
$ws.Range("J83").Value = 2876.875
$ws.Range("K83").Value = 42612.5
$ws.Range("L83").Value = 14384.375
$ws.Range("M83").Value = -37620.5
$ws.Range("N83").Value = -24368.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 441.5625
$ws.Range("I22").Value = 414.75
$ws.Range("J22").Value = 522
$ws.Range("K22").Value = 414.75
$ws.Range("L22").Value = 522
$ws.Range("M22").Value = -119.75
$ws.Range("N22").Value = -1112

$ws.Range("H27").Value = 441.5625
$ws.Range("I27").Value = 414.75
$ws.Range("J27").Value = 522
$ws.Range("K27").Value = 414.75
$ws.Range("L27").Value = 522
$ws.Range("M27").Value = -307.75
$ws.Range("N27").Value = -736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 14429
$ws.Range("J123").Value = 14429
$ws.Range("L123").Value = 14429
$ws.Range("N123").Value = -24229

Write-Output "Updated 163 cells across 8 sheets."